# Apply the diff:
#  - slide 5, shape 2 ("ash  and  proanthocyanins..."): reposition (offset only)
#  - slide 6, shape 1 ("So the final array..."): reposition/resize slightly
#  - slide 6, shape 2 (the big bracketed tuple list): rename, reposition/resize,
#    enable word-wrap, and reflow the single paragraph into 13 separate paragraphs
#  - slide 6, shape 3 ("You can see that..."): reposition/resize slightly
#  - slide 6, shape 4 ("Ash (2) Proanthocyanins..."): reposition (offset only)
#
# Note: Shape.Left/Top/Width/Height are single-precision (Single) in the
# PowerPoint object model, and this host floors (rather than rounds) the
# point->EMU conversion, so the literals below were chosen so that
# floor(Single(points) * 12700) lands exactly on the target EMU values from
# the target OOXML.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 5 / shape 2 ("ash  and  proanthocyanins...")
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$sh = $s5.Shapes.Item(2)
$sh.Left = 283.50003050045655
$sh.Top  = 181.27217108413055

# ---------------------------------------------------------------------------
# Slide 6
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)

# shape 1: "So the final array for incremental removal is as follows:"
$sh1 = $s6.Shapes.Item(1)
$sh1.Left   = 100.000038199879
$sh1.Top    = 24.446889893724624
$sh1.Width  = 824.0000610992919
$sh1.Height = 89.00011444003663

# shape 2: the long bracketed tuple list -- rename + reflow into paragraphs
$sh2 = $s6.Shapes.Item(2)
$sh2.Name = "[(2, " + [char]0x2018 + "ash" + [char]0x2019 + ")," + [char]0x2026
$sh2.TextFrame.WordWrap = -1
$sh2.TextFrame.TextRange.Text = "[(2, " + [char]0x2018 + "ash" + [char]0x2019 + "),`r(2, " + [char]0x2018 + "proanthocyanins'),`r(2, " + [char]0x2018 + "color_intensity'),`r (1, 'alcalinity_of_ash'), `r(1, " + [char]0x2018 + "od280/od315_of_diluted_wines" + [char]0x2019 + "),`r (0, " + [char]0x2018 + "alcohol'),`r (0, " + [char]0x2018 + "malic_acid'),`r (0, " + [char]0x2018 + "magnesium'),`r (0, " + [char]0x2018 + "total_phenols'),`r (0, " + [char]0x2018 + "flavanoids" + [char]0x2019 + "),`r (0, " + [char]0x2018 + "nonflavanoid_phenols'),`r (0, 'hue'), `r(0, 'proline')]"
$sh2.Left   = 216.58396910775338
$sh2.Top    = 112.62263492507596
$sh2.Width  = 590.8321533835555
$sh2.Height = 413.30398560833464

# shape 3: "You can see that the 5 first features are represented with their weight"
$sh3 = $s6.Shapes.Item(3)
$sh3.Left   = 125.19184882350108
$sh3.Top    = 517.3833618859082
$sh3.Width  = 798.8640747473241
$sh3.Height = 36.30397800789787

# shape 4: "Ash (2) Proanthocyanins (2)..."
$sh4 = $s6.Shapes.Item(4)
$sh4.Top = 566.5074462741006
